$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting rows 106:160 down to 107:161
$ws.Rows("106:106").Insert()

# Populate new row 106 with the new data row
$ws.Range("A106").Value = 4
$ws.Range("B106").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C106").Value = "Los Lagos"
$ws.Range("D106").Value = 44572
$ws.Range("D106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E106").Value = 10
$ws.Range("F106").Value = 100112039
$ws.Range("G106").Value = "Ciboulette"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 240
$ws.Range("K106").Value = 3000
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = 3000
$ws.Range("N106").Value = "$/docena de atados"
$ws.Range("O106").Value = "Región Metropolitana"
$ws.Range("P106").Value = 1000
$ws.Range("Q106").Value = 3
$ws.Range("R106").Value = "Hortaliza"
